# Updated cryptos list on Fri Mar  1 01:40:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value. Column D sometimes holds
# plain decimal-looking numbers (e.g. "402.80") that must stay TEXT (they
# are thousands-dotted price strings elsewhere, e.g. "61.358.44"), so for
# those we force a text format before the write and restore the default
# style afterwards so the cell format matches the original workbook.
$updates = [ordered]@{
    "D2" = "61.358.44"
    "E2" = "  -0.11%  "
    "D3" = "3.376.34"
    "E3" = "  -1.40%  "
    "E4" = "  +0.25%  "
    "D5" = "402.80"
    "E5" = "  -2.84%  "
    "D6" = "131.93"
    "E6" = "  +8.42%  "
    "D7" = "0.590"
    "E7" = "  +2.00%  "
    "E8" = "  +0.06%  "
    "D9" = "0.667"
    "E9" = "  +3.78%  "
    "E10" = "  +5.42%  "
    "D11" = "41.59"
    "E11" = "  +0.78%  "
    "E12" = "  -0.93%  "
    "D13" = "3.916.11"
    "E13" = "  -1.15%  "
    "D14" = "8.34"
    "E14" = "  -1.27%  "
    "D15" = "19.55"
    "E15" = "  -0.32%  "
    "D16" = "3.366.16"
    "E16" = "  -1.49%  "
    "D17" = "61.376.67"
    "E17" = "  +0.16%  "
    "D18" = "1.02"
    "E18" = "  -0.91%  "
    "D19" = "11.15"
    "E19" = "  +1.96%  "
    "D20" = "0.0000126"
    "E20" = "  +7.76%  "
    "D21" = "3.20"
    "E21" = "  -4.96%  "
    "D22" = "83.01"
    "E22" = "  +9.21%  "
    "D23" = "12.77"
    "E23" = "  -2.04%  "
    "D24" = "303.84"
    "E24" = "  +1.78%  "
    "D25" = "3.12"
    "E25" = "  -0.07%  "
    "D26" = "4.78"
    "E26" = "  +12.00%  "
    "D27" = "8.40"
    "E27" = "  +7.96%  "
    "D28" = "29.20"
    "E28" = "  -5.17%  "
    "D29" = "7.75"
    "E29" = "  -3.21%  "
    "E30" = "  +0.44%  "
    "D31" = "0.116"
    "E31" = "  +1.46%  "
    "E32" = "  +0.17%  "
    "D33" = "11.28"
    "E33" = "  -1.41%  "
    "D34" = "41.29"
    "E34" = "  -3.60%  "
    "D35" = "2.49"
    "E35" = "  -0.81%  "
    "D36" = "0.0477"
    "E36" = "  -1.31%  "
    "D37" = "51.83"
    "E37" = "  -0.46%  "
    "D38" = "1.00"
    "E38" = "  +0.31%  "
    "D39" = "3.40"
    "E39" = "  -3.05%  "
    "D40" = "2.93"
    "E40" = "  -3.40%  "
    "D41" = "137.47"
    "E41" = "  +1.91%  "
    "D42" = "1.98"
    "E42" = "  +0.29%  "
    "D43" = "0.124"
    "E43" = "  +1.57%  "
    "D44" = "0.289"
    "E44" = "  +1.49%  "
    "D45" = "3.96"
    "E45" = "  -0.08%  "
    "D46" = "16.71"
    "E46" = "  -3.65%  "
    "D47" = "2.23"
    "E47" = "  +1.66%  "
    "D48" = "21.34"
    "E48" = "  -4.94%  "
    "B49" = "Maker"
    "C49" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D49" = "2.112.63"
    "E49" = "  -4.28%  "
    "B50" = "ApeXProtocol"
    "C50" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D50" = "2.29"
    "E50" = "  -4.44%  "
    "B51" = "ThetaToken"
    "C51" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "D51" = "1.86"
    "E51" = "  -0.72%  "
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    $isPriceColumn = $ref.StartsWith("D")
    $looksNumeric = $false
    if ($isPriceColumn) {
        $looksNumeric = $value -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    }
    if ($looksNumeric) {
        # Force text storage so e.g. "402.80" / "1.00" is not coerced into
        # the number 402.8 / 1, then restore the default cell style so the
        # saved format matches the rest of the (unstyled) data cells.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
